$wb = $excel.ActiveWorkbook

# --- Sheet: Full results ---
$ws1 = $wb.Worksheets.Item("Full results")

# Rename shared string value "income_hh" -> "income" by updating all cells that reference it
$ws1.Range("A8").Value = "income"
$ws1.Range("A9").Value = "income"
$ws1.Range("A10").Value = "income"

# Updated numeric results (tiny precision differences from recalculation)
$ws1.Range("L2").Value = 0.0140135097403642
$ws1.Range("M2").Value = 0.0440429494849118
$ws1.Range("N2").Value = 0.320043830435395
$ws1.Range("H4").Value = 0.602270479269944
$ws1.Range("I4").Value = 0.285664744672762
$ws1.Range("O4").Value = 0.397864820589618
$ws1.Range("L5").Value = 0.0258778380577954
$ws1.Range("M5").Value = 0.0575515385688163
$ws1.Range("N5").Value = 0.150539331758559
$ws1.Range("H7").Value = 0.819804882719041
$ws1.Range("I7").Value = 0.125102165652733
$ws1.Range("O7").Value = 0.180258984609072
$ws1.Range("L8").Value = 0.0418768880182208
$ws1.Range("M8").Value = 0.0497546070209868
$ws1.Range("N8").Value = 0.197026389869398
$ws1.Range("H10").Value = 0.785914068190657
$ws1.Range("I10").Value = 0.129262717270778
$ws1.Range("O10").Value = 0.214167972619614
$ws1.Range("L11").Value = 0.0283411036041199
$ws1.Range("M11").Value = 0.0328616107336478
$ws1.Range("N11").Value = 0.253728365699651
$ws1.Range("H13").Value = 0.735913239232575
$ws1.Range("I13").Value = 0.220052021734735
$ws1.Range("O13").Value = 0.264194850096421

# --- Sheet: For plotting ---
$ws2 = $wb.Worksheets.Item("For plotting")

# Rename shared string value "income_hh" -> "income"
$ws2.Range("E8").Value = "income"
$ws2.Range("E9").Value = "income"
$ws2.Range("E10").Value = "income"

# Updated numeric results (tiny precision differences from recalculation)
$ws2.Range("B2").Value = 0.270173855721891
$ws2.Range("C2").Value = 0.369913805148899
$ws2.Range("D2").Value = 0.320043830435395
$ws2.Range("B3").Value = 0.356556663617972
$ws2.Range("C3").Value = 0.439172977561264
$ws2.Range("D3").Value = 0.397864820589618
$ws2.Range("B4").Value = 0.308862869302389
$ws2.Range("C4").Value = 0.398780872907023
$ws2.Range("B5").Value = 0.0869753789565366
$ws2.Range("C5").Value = 0.214103284560582
$ws2.Range("D5").Value = 0.150539331758559
$ws2.Range("B6").Value = 0.125933173247459
$ws2.Range("C6").Value = 0.234584795970684
$ws2.Range("D6").Value = 0.180258984609072
$ws2.Range("B7").Value = 0.0692263519047685
$ws2.Range("C7").Value = 0.176188540175742
$ws2.Range("B8").Value = 0.0771328980096298
$ws2.Range("C8").Value = 0.316919881729166
$ws2.Range("D8").Value = 0.197026389869398
$ws2.Range("B9").Value = 0.100892028059918
$ws2.Range("C9").Value = 0.32744391717931
$ws2.Range("D9").Value = 0.214167972619614
$ws2.Range("B10").Value = 0.0341328586480976
$ws2.Range("C10").Value = 0.294693872549156
$ws2.Range("B11").Value = 0.0867093433616834
$ws2.Range("C11").Value = 0.420747388037619
$ws2.Range("D11").Value = 0.253728365699651
$ws2.Range("B12").Value = 0.110305346420244
$ws2.Range("C12").Value = 0.418084353772599
$ws2.Range("D12").Value = 0.264194850096421
$ws2.Range("B13").Value = 0.0556392793389844
$ws2.Range("C13").Value = 0.407027199386563
$ws2.Range("B14").Value = 0.0759450637572911
$ws2.Range("C14").Value = 0.213212731355565
$ws2.Range("B15").Value = 0.103132067181238
$ws2.Range("C15").Value = 0.229616782012906
$ws2.Range("B16").Value = 0.0725894010190001
$ws2.Range("C16").Value = 0.206436618702007
